$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$src = $ws1.Range("B1")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "Design Info"
$newSheet.Activate()
$newSheet.Range("B9").Value2 = "NOTE"
$src.Copy()
$newSheet.Range("B9").PasteSpecial(-4122)
$newSheet.Range("B9").Font.Bold = $true
Write-Output "done"
